$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns C (record_atd) and D (average_simulation_TD), rows 2-18
$values = @{
    2  = @(9, 10.5)
    3  = @(7, 9)
    4  = @(28, 30)
    5  = @(1, 2)
    6  = @(51, 45)
    7  = @(29, 32)
    8  = @(56, 55)
    9  = @(50, 45)
    10 = @(12, 10)
    11 = @(2, 1.5)
    12 = @(22, 22.5)
    13 = @(55, 61.5)
    14 = @(5, 5)
    15 = @(8, 9)
    16 = @(3, 3.5)
    17 = @(145, 151)
    18 = @(14, 13.5)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 3).Value = $pair[0]
    $ws.Cells.Item($row, 4).Value = $pair[1]
}

# Row 19 is the averages row - only column C has a value (no D19 cell)
$ws.Cells.Item(19, 3).Value = 29.23529411764706
